$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected range to Text format first so that Excel stores the
# updated values as text (matching the original inline-string cells) rather
# than auto-converting numeric-looking strings into numbers.
$rng = $ws.Range("B2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "56.417.47"
$ws.Range("E2").Value = "  +9.52%  "
$ws.Range("D3").Value = "3.235.60"
$ws.Range("E3").Value = "  +4.17%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "399.14"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").Value = "111.27"
$ws.Range("E6").Value = "  +6.92%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  +6.06%  "
$ws.Range("D10").Value = "39.53"
$ws.Range("E10").Value = "  +6.16%  "
$ws.Range("E11").Value = "  +5.65%  "
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "3.744.46"
$ws.Range("E13").Value = "  +4.23%  "
$ws.Range("D14").Value = "8.13"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("D15").Value = "19.12"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "1.06"
$ws.Range("E16").Value = "  +5.55%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.209.46"
$ws.Range("E17").Value = "  +3.72%  "
$ws.Range("D18").Value = "10.70"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "56.200.67"
$ws.Range("E19").Value = "  +8.95%  "
$ws.Range("D20").Value = "3.33"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("E21").Value = "  +5.75%  "
$ws.Range("D22").Value = "13.09"
$ws.Range("E22").Value = "  +4.36%  "
$ws.Range("D23").Value = "305.70"
$ws.Range("E23").Value = "  +14.29%  "
$ws.Range("E24").Value = "  +7.65%  "
$ws.Range("D25").Value = "3.24"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "28.32"
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("D28").Value = "7.48"
$ws.Range("E28").Value = "  +3.17%  "
$ws.Range("E29").Value = "  +4.26%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +3.84%  "
$ws.Range("D32").Value = "11.21"
$ws.Range("E32").Value = "  +7.28%  "
$ws.Range("D33").Value = "0.0494"
$ws.Range("E33").Value = "  +3.52%  "
$ws.Range("D34").Value = "36.48"
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").Value = "51.39"
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("D37").Value = "3.15"
$ws.Range("E37").Value = "  +25.08%  "
$ws.Range("D38").Value = "3.53"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "135.22"
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("D42").Value = "4.03"
$ws.Range("E42").Value = "  +6.52%  "
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("D44").Value = "0.119"
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("D46").Value = "22.34"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  +47.16%  "
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").Value = "2.47"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").Value = "2.139.38"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("D51").Value = "0.0364"
$ws.Range("E51").Value = "  +8.12%  "

# Restore the default style so we do not leave a stray number-format
# applied to the cells (keeps styles identical to the original workbook).
$rng.Style = "Normal"
